# P2p implementation changes and completed status
#
# The login-data sheet's user list advances from the "auttestt_NN" series
# to a new "Sanuserr_NN" series, shifted up by one entry, and extended with
# six additional rows (new users Sanuserr_25 .. Sanuserr_30), for a total of
# 18 data rows (rows 2-19). Column B keeps the constant "password" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2) starts at Sanuserr_13 and increments through Sanuserr_30 at row 19.
$startNum = 13
$firstDataRow = 2
$lastDataRow = 19

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $num = $startNum + ($row - $firstDataRow)
    $ws.Cells.Item($row, 1).Value = "Sanuserr_$num@mailinator.com"
    $ws.Cells.Item($row, 2).Value = "password"
}

# Reflect the completed paste/entry range in the sheet's active selection.
$ws.Range("A2:B21").Select() | Out-Null
